# Agregado análisis avanzado Planes
#
# The PLANES worksheet's feature-matrix table is restructured:
#   - The plan names in row 1 are changed from ALL-CAPS to Title-Case
#     (FREE -> Free, BASICO -> Basico, PRO -> Pro, PREMIUM -> Premium).
#   - The feature rows are renamed/reordered:
#       VER_KPIS        -> TARJETAS_KPI
#       CONCLUSIONES    -> ANALISIS_LENGUAJE
#       EXPORTAR_EXCEL  -> ANALISIS_AVANZADO
#       EXPORTAR_PDF    stays EXPORTAR_PDF
#   - The old last row (USOS_MAXIMOS / 3,9999,9999,9999) is removed.
#   - Column widths are widened to fit the new content.
#   - The PLANES tab becomes the active/selected sheet.
#   - Page setup (paper size / orientation) is defined for PLANES.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PLANES")

# --- Remove the trailing USOS_MAXIMOS row (old row 6) ---
$ws.Rows.Item(6).Delete()

# --- Row 1: plan names, Title-Case instead of ALL CAPS ---
$ws.Range("B1").Value = "Free"
$ws.Range("C1").Value = "Basico"
$ws.Range("D1").Value = "Pro"
$ws.Range("E1").Value = "Premium"

# --- Column A: renamed functionality keys (row 5 "EXPORTAR_PDF" unchanged) ---
$ws.Range("A2").Value = "TARJETAS_KPI"
$ws.Range("A3").Value = "ANALISIS_LENGUAJE"
$ws.Range("A4").Value = "ANALISIS_AVANZADO"

# --- Widen columns to fit the new headers/content ---
$ws.Columns.Item(1).ColumnWidth = 21.28
$ws.Columns.Item(2).ColumnWidth = 17.61
$ws.Columns.Item(3).ColumnWidth = 18.94
$ws.Columns.Item(4).ColumnWidth = 19.61
$ws.Columns.Item(5).ColumnWidth = 19.17
$ws.Columns.Item(6).ColumnWidth = 18.17

# --- Page setup for PLANES ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Make PLANES the active/selected sheet ---
$ws.Activate()
